$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to reset style
# after forcing text entry on cells whose new value would otherwise be
# auto-parsed by Excel as a number (e.g. "1.002").
$defaultStyle = $ws.Range("A100").Style

$ws.Range("D2").Value = "27.489.76"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").Value = "1.748.04"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").Value = "'324.30"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").Value = "'0.4460"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +4.41%  "

$ws.Range("D8").Value = "'0.3587"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").Value = "'0.07505"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +0.36%  "

$ws.Range("D10").Value = "'42.01"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -5.44%  "

$ws.Range("D11").Value = "'1.091"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -1.48%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").Value = "'20.74"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -3.90%  "

$ws.Range("D14").Value = "'6.010"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -1.95%  "

$ws.Range("D15").Value = "'7.105"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -2.94%  "

$ws.Range("D16").Value = "1.750.44"
$ws.Range("E16").Value = "  -2.17%  "

$ws.Range("D17").Value = "'93.10"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "'0.06404"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").Value = "'16.79"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "'5.801"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").Value = "27.549.25"
$ws.Range("E23").Value = "  -1.36%  "

$ws.Range("D24").Value = "'11.16"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("D25").Value = "'2.095"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -2.92%  "

$ws.Range("D26").Value = "'162.89"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +1.82%  "

$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("D28").Value = "1.950.16"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("D29").Value = "'2.081"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  -3.95%  "

$ws.Range("D30").Value = "'125.70"
$ws.Range("D30").Style = $defaultStyle

$ws.Range("E31").Value = "  -8.29%  "

$ws.Range("D32").Value = "'3.664"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +4.60%  "

$ws.Range("D33").Value = "'0.09055"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("D34").Value = "'5.518"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -2.99%  "

$ws.Range("E35").Value = "  -5.92%  "

$ws.Range("D36").Value = "'0.02278"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -1.64%  "

$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").Value = "'0.6353"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("D39").Value = "'0.05989"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").Value = "'4.935"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("D41").Value = "'1.200"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("D42").Value = "'1.381"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  -0.50%  "

$ws.Range("D43").Value = "'7.741"
$ws.Range("D43").Style = $defaultStyle

$ws.Range("D44").Value = "'13.14"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -3.68%  "

$ws.Range("D45").Value = "'3.714"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("D46").Value = "'0.5868"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -1.70%  "

$ws.Range("D47").Value = "'121.62"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -2.10%  "

$ws.Range("D48").Value = "'1.945"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("D49").Value = "'1.140"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("D50").Value = "'0.06835"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("E51").Value = "  -3.01%  "
